$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2036
$ws.Range("F3").Value = 657
$ws.Range("F4").Value = 1303
$ws.Range("F7").Value = 154
$ws.Range("F8").Value = 371
$ws.Range("F9").Value = 159
$ws.Range("F10").Value = 119
$ws.Range("F11").Value = 941
$ws.Range("F12").Value = 298
$ws.Range("F13").Value = 157
$ws.Range("F14").Value = 39
$ws.Range("F17").Value = 318
$ws.Range("F18").Value = 733
$ws.Range("F19").Value = 112
$ws.Range("F20").Value = 693
$ws.Range("F21").Value = 235
$ws.Range("F23").Value = 947
$ws.Range("F24").Value = 409
$ws.Range("F25").Value = 225
$ws.Range("F26").Value = 73
$ws.Range("F27").Value = 336
$ws.Range("F30").Value = 443

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 346
$ws.Range("F5").Value = 28
$ws.Range("F6").Value = 36
$ws.Range("F7").Value = 274

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 339

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 339
$ws.Range("F3").Value = 2036
$ws.Range("F4").Value = 657
$ws.Range("F5").Value = 1303
$ws.Range("F9").Value = 154
$ws.Range("F10").Value = 371
$ws.Range("F11").Value = 159
$ws.Range("F12").Value = 119
$ws.Range("F13").Value = 941
$ws.Range("F14").Value = 298
$ws.Range("F15").Value = 157
$ws.Range("F17").Value = 39
$ws.Range("F18").Value = 346
$ws.Range("F20").Value = 28
$ws.Range("F22").Value = 36
$ws.Range("F23").Value = 274
$ws.Range("F24").Value = 318
$ws.Range("F25").Value = 733
$ws.Range("F26").Value = 112
$ws.Range("F27").Value = 693
$ws.Range("F28").Value = 235
$ws.Range("F30").Value = 947
$ws.Range("F31").Value = 409
$ws.Range("F34").Value = 225
$ws.Range("F35").Value = 73
$ws.Range("F36").Value = 336
$ws.Range("F42").Value = 443

